$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artilharia")

# Add a "Nome" (player name) column in F, mirroring column B.
$ws.Range("F1").Value = "Nome"
$ws.Range("F3").Value = "Flaco López"
$ws.Range("F4").Value = "Dellatorre"
$ws.Range("F5").Value = "Raphael Veiga"
$ws.Range("F6").Value = "Raphael Veiga"
$ws.Range("F7").Value = "Yuri Alberto"
$ws.Range("F8").Value = "Eduardo Sasha"
$ws.Range("F9").Value = "Jeferson Jeh"
$ws.Range("F10").Value = "Jenison"

# Copy column B's formatting (thin border, same style already used by the
# rest of the table) onto the new column F instead of defining a new style.
$ws.Range("B1:B10").Copy() | Out-Null
$ws.Range("F1:F10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Columns.Item(6).ColumnWidth = 12.3

# Normalize E1's style (drop the stray fill flag) to match the rest of row 1.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Make "Artilharia" the active sheet/tab and set its selection, matching the
# saved view state in the workbook.
$ws.Activate()
$ws.Range("F15").Select() | Out-Null
